# AutoCommit_17 июня 2024 г. 20:21:12_SibNout2023
# "Расставить двойки" - fill in missing (0) grades with 2 ("неуд"),
# and flag row 24 (a student with no grades at all) with a purple highlight.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Make-BGR($r, $g, $b) {
    return ($b * 65536) + ($g * 256) + $r
}

# Purple highlight colour used for row 24 (RGB 70,30,A0 -> FF7030A0 in the sheet).
$purple = Make-BGR 0x70 0x30 0xA0

# --- Row 3: header-ish row that was blank, now filled with 1s (C3:G3) ---
$ws.Range("C3:G3").Value = 1

# --- Replace every "0" placeholder grade with "2" (an F grade), row by row ---
# (the engine only honours the first area of a multi-area Range write, so each
# contiguous block is set individually instead of as one big union range)
$zeroAreas = @(
    "G4:J4",
    "G7",
    "C8",
    "G10:J10",
    "G13:J13",
    "G15:J15",
    "I17:J17",
    "G18:J18",
    "G21:J21",
    "G25:J25",
    "C30",
    "H30:J30",
    "J32",
    "J33",
    "J34",
    "J35",
    "J36",
    "J37",
    "J38",
    "J39",
    "J40",
    "J41",
    "J42",
    "J43",
    "J44",
    "J45",
    "J46",
    "J47",
    "J48",
    "J49",
    "J50"
)
foreach ($area in $zeroAreas) {
    $ws.Range($area).Value = 2
}

# --- Row 24: every grade was 0 (student never showed up) -> set to 2 and
#     highlight the row in purple, adding a purple spacer cell in K24 ---
$ws.Range("C24:J24").Value = 2
$ws.Range("D24:J24").Interior.Color = $purple
$ws.Range("K24").Interior.Color = $purple

# --- Conditional formatting: carve the now explicitly-coloured D24:J24 out of
#     the colour-scale range applied to the grade columns ---
$cf = $ws.Range("C4:J50,M4:M51").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("C4:J23,M4:M51,C25:J50,C24"))

# --- Selection: active cell moves from B3 to G4 ---
$ws.Range("G4").Select()
